# chore(runtime): publish files + archive (2025-12-07 11:07:34)
#
# Updates the KHL stats export:
#  - Matches_SOG: append two newly finished matches (rows 364-365)
#  - Shots_HA / Shots_Summary / Meta_ext: roll the as_of_utc snapshot
#    timestamp forward and refresh the per-team shots-on-goal aggregates
#    that moved because of the two new matches.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Matches_SOG - append the two new matches
# ---------------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

# uid values are stored as text (matches the existing rows, which are all
# inline strings) - a leading apostrophe forces Excel to keep the
# numeric-looking uid as text instead of coercing it to a number.
$wsMatches.Range("A364").Value = "'897857"
$wsMatches.Range("B364").Value = "2025-12-06T10:00:00"
$wsMatches.Range("C364").Value = "Адмирал"
$wsMatches.Range("D364").Value = "Локомотив"
$wsMatches.Range("E364").Value = 29
$wsMatches.Range("F364").Value = 25
$wsMatches.Range("G364").Value = "khl_text"

$wsMatches.Range("A365").Value = "'897858"
$wsMatches.Range("B365").Value = "2025-12-06T10:00:00"
$wsMatches.Range("C365").Value = "Амур"
$wsMatches.Range("D365").Value = "СКА"
$wsMatches.Range("E365").Value = 42
$wsMatches.Range("F365").Value = 37
$wsMatches.Range("G365").Value = "khl_text"

# ---------------------------------------------------------------------
# Sheet: Shots_HA - refresh as_of_utc for every team row, plus the
# home/away shots-on-goal totals for the four teams that played.
# ---------------------------------------------------------------------
$wsHA = $wb.Worksheets.Item("Shots_HA")
$newAsOf = "2025-12-06T10:00:00Z"

for ($r = 2; $r -le 23; $r++) {
    $wsHA.Range("D$r").Value = $newAsOf
}

# Row 4: Адмирал (home in the new match)
$wsHA.Range("E4").Value = 16
$wsHA.Range("G4").Value = 584
$wsHA.Range("H4").Value = 433
$wsHA.Range("I4").Value = 36.5
$wsHA.Range("J4").Value = 27.1

# Row 6: Амур (home in the new match)
$wsHA.Range("E6").Value = 18
$wsHA.Range("G6").Value = 548
$wsHA.Range("H6").Value = 618
$wsHA.Range("I6").Value = 30.4
$wsHA.Range("J6").Value = 34.3

# Row 12: Локомотив (away in the new match)
$wsHA.Range("F12").Value = 19
$wsHA.Range("K12").Value = 576
$wsHA.Range("L12").Value = 470
$wsHA.Range("M12").Value = 30.3
$wsHA.Range("N12").Value = 24.7

# Row 15: СКА (away in the new match)
$wsHA.Range("F15").Value = 17
$wsHA.Range("K15").Value = 513
$wsHA.Range("L15").Value = 583
$wsHA.Range("M15").Value = 30.2
$wsHA.Range("N15").Value = 34.3

# ---------------------------------------------------------------------
# Sheet: Shots_Summary - refresh as_of_utc for every team row, plus the
# combined totals for the same four teams.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Shots_Summary")

for ($r = 2; $r -le 23; $r++) {
    $wsSummary.Range("D$r").Value = $newAsOf
}

# Row 4: Адмирал
$wsSummary.Range("E4").Value = 32
$wsSummary.Range("F4").Value = 1081
$wsSummary.Range("G4").Value = 877
$wsSummary.Range("H4").Value = 33.8
$wsSummary.Range("I4").Value = 27.4

# Row 6: Амур
$wsSummary.Range("E6").Value = 34
$wsSummary.Range("F6").Value = 986
$wsSummary.Range("G6").Value = 1232
$wsSummary.Range("H6").Value = 29

# Row 12: Локомотив
$wsSummary.Range("E12").Value = 36
$wsSummary.Range("F12").Value = 1148
$wsSummary.Range("G12").Value = 920
$wsSummary.Range("H12").Value = 31.9
$wsSummary.Range("I12").Value = 25.6

# Row 15: СКА
$wsSummary.Range("E15").Value = 32
$wsSummary.Range("F15").Value = 1001
$wsSummary.Range("G15").Value = 1079
$wsSummary.Range("H15").Value = 31.3
$wsSummary.Range("I15").Value = 33.7

# ---------------------------------------------------------------------
# Sheet: Meta_ext - bump the snapshot timestamp and build_version
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Range("B2").Value = $newAsOf
$wsMeta.Range("D2").Value = 38
